# Apply updated dSF (column F) values for the re-pulled / recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 2
    13 = 9
    16 = 3
    20 = 1
    26 = 0
    29 = -3
    31 = -8
    33 = -3
    34 = 0
    35 = 3
    36 = -5
    40 = 6
    44 = -10
    46 = -8
    48 = -9
    51 = -6
    54 = -6
    58 = -6
    60 = -6
    61 = 0
    62 = -2
    68 = -7
    69 = -10
    70 = 0
    71 = -2
    72 = -1
    73 = 0
    76 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
